$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 3.1.0 -> 3.2.0
$ws.Range("B3").Value = "3.2.0"

# Title: Dk Core SNOMED CT Basic Observation -> DK Core SNOMED CT Basic Observation
$ws.Range("B5").Value = "DK Core SNOMED CT Basic Observation"

# Date: 2024-01-08T21:51:35+01:00 -> 2024-05-06T15:28:33+02:00
$ws.Range("B8").Value = "2024-05-06T15:28:33+02:00"

# Contact: No display for ContactDetail -> HL7 Denmark (http://www.hl7.dk, jenskristianvilladsen@gmail.com)
$ws.Range("B10").Value = "HL7 Denmark (http://www.hl7.dk, jenskristianvilladsen@gmail.com)"
